$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("F22", 4),
    @("G22", 122.96),
    @("F26", 24),
    @("G26", 1217.28),
    @("F27", 40),
    @("G27", 1639.6),
    @("B40", 69031.92999999999),
    @("F61", 56),
    @("G61", 2497.6),
    @("F62", 146),
    @("G62", 8140.96),
    @("F64", 52),
    @("G64", 1836.64),
    @("F65", 100),
    @("G65", 7793),
    @("B73", 250476.55),
    @("F80", 2),
    @("G80", 1143.94),
    @("B100", 56928.48),
    @("F102", 1),
    @("G102", 2320.17),
    @("B104", 30338.01),
    @("F134", 103),
    @("G134", 4814.22),
    @("B145", 86819.13),
    @("F185", 45),
    @("G185", 6001.2),
    @("B189", 42621.51),
    @("F197", 0),
    @("G197", 0),
    @("B199", -6196.5),
    @("F221", 73),
    @("G221", 5738.53),
    @("F231", 50),
    @("G231", 4477),
    @("B233", 63255),
    @("F233", 73),
    @("G233", 5986),
    @("B234", 57004),
    @("F234", 0),
    @("G234", 0),
    @("F236", 55),
    @("G236", 2380.4),
    @("B238", 57552),
    @("E238", 136.86),
    @("F238", -5),
    @("G238", -603.45),
    @("B239", 64329),
    @("E239", 128.32),
    @("F239", 0),
    @("G239", 0),
    @("F242", 12),
    @("G242", 1068.72),
    @("F245", 32),
    @("G245", 1916.48),
    @("B247", 82578.89),
    @("F272", 153),
    @("G272", 8324.73),
    @("F273", 12),
    @("G273", 1224.12),
    @("B280", 95950.56),
    @("F284", 1619),
    @("G284", 29951.5),
    @("B291", 47255.43),
    @("F341", 36),
    @("G341", 4124.16),
    @("B371", 131224.15),
    @("F404", 138),
    @("G404", 23643.54),
    @("B408", 24514.48),
    @("F433", 9),
    @("G433", 3073.23),
    @("B450", 88696.97),
    @("F455", 311),
    @("G455", 43723.49),
    @("B457", 98202.97),
    @("F498", 283),
    @("G498", 27804.75),
    @("F501", 75),
    @("G501", 2573.25),
    @("F502", 98),
    @("G502", 3972.92),
    @("B518", 192775.51),
    @("B555", 64922),
    @("E555", 20.98),
    @("F555", 0),
    @("G555", 0),
    @("B556", 45706),
    @("E556", 23.58),
    @("F556", -207),
    @("G556", -4084.11),
    @("B568", 53595),
    @("E568", 17.61),
    @("F568", -338),
    @("G568", -4978.74),
    @("B569", 65067),
    @("E569", 15.65),
    @("F569", 0),
    @("G569", 0),
    @("B658", 60025),
    @("E658", 37.22),
    @("F658", -98),
    @("G658", -3217.34),
    @("B659", 64833),
    @("E659", 34.9),
    @("F659", 88),
    @("G659", 2889.04),
    @("F660", 132),
    @("G660", 4333.56),
    @("F663", 75),
    @("G663", 7387.5),
    @("F667", 163),
    @("G667", 13701.78),
    @("B668", 60022),
    @("E668", 37.22),
    @("F668", -113),
    @("G668", -3709.79),
    @("B669", 64830),
    @("E669", 34.9),
    @("F669", 88),
    @("G669", 2889.04),
    @("B670", 42243.98),
    @("F688", 185),
    @("G688", 15837.85),
    @("B691", 151611.79),
    @("F700", 254),
    @("G700", 8409.940000000001),
    @("F703", 166),
    @("G703", 7264.16),
    @("B704", 33717.46),
    @("F734", 277),
    @("G734", 33752.45),
    @("F739", 61),
    @("G739", 4903.18),
    @("B741", 41046.25),
    @("F804", 2),
    @("G804", 327.78),
    @("F806", 77),
    @("G806", 8378.370000000001),
    @("F810", 265),
    @("G810", 39858.65),
    @("F818", 55),
    @("G818", 2646.6),
    @("F822", 5),
    @("G822", 165.3),
    @("F827", 94),
    @("G827", 48366.76),
    @("F830", 375),
    @("G830", 13811.25),
    @("F831", 45),
    @("G831", 2124.45),
    @("B837", 262156.8),
    @("F876", 75),
    @("G876", 6023.25),
    @("B882", 18567.6),
    @("F887", 54),
    @("G887", 1632.42),
    @("B894", 244869.94),
    @("F905", 6),
    @("G905", 245.94),
    @("F907", 44),
    @("G907", 5515.4),
    @("B910", 14787.16),
    @("B937", 3611924.28),
    @("B938", 3611924.28)
)

foreach ($u in $updates) {
    $ws.Range($u[0]).Value = $u[1]
}
